{"js": "// Replacement pairs extracted from the diff, in document order (old text -> new text).\n// All old texts are unique within the document, so this is unambiguous.\nconst replacements = [\n    [\"96-9=\", \"3+61=\"],\n    [\"55-41=\", \"75-30=\"],\n    [\"26+40=\", \"33+40=\"],\n    [\"88-2=\", \"63-12=\"],\n    [\"87-46=\", \"19+20=\"],\n    [\"54-3=\", \"51+44=\"],\n    [\"41+26=\", \"81-20=\"],\n    [\"82-3=\", \"46-29=\"],\n    [\"32+42=\", \"82-19=\"],\n    [\"84+1=\", \"0+52=\"],\n    [\"17+1=\", \"77-50=\"],\n    [\"68-22=\", \"42+44=\"],\n    [\"25+17=\", \"68-7=\"],\n    [\"95-86=\", \"9+64=\"],\n    [\"72-65=\", \"3+71=\"],\n    [\"29+32=\", \"45-44=\"],\n    [\"27+57=\", \"47+44=\"],\n    [\"40+27=\", \"87-64=\"],\n    [\"32+38=\", \"53+5=\"],\n    [\"86-48=\", \"8+14=\"],\n    [\"25+73=\", \"59+22=\"],\n    [\"73-49=\", \"42-32=\"],\n    [\"11+19=\", \"90-80=\"],\n    [\"89-19=\", \"54+40=\"],\n    [\"37-22=\", \"68-65=\"],\n    [\"52-24=\", \"17+56=\"],\n    [\"80-4=\", \"2+71=\"],\n    [\"38+50=\", \"22-11=\"],\n    [\"13+69=\", \"46+8=\"],\n    [\"16+78=\", \"33+34=\"],\n    [\"36+58=\", \"14+14=\"],\n    [\"30+55=\", \"13+9=\"],\n    [\"39-30=\", \"91-2=\"],\n    [\"23-17=\", \"88-0=\"],\n    [\"12+42=\", \"4+16=\"],\n    [\"5+0=\", \"11+67=\"],\n    [\"32+9=\", \"79-33=\"],\n    [\"4+49=\", \"69+27=\"],\n    [\"40-25=\", \"82-29=\"],\n    [\"92-0=\", \"28+9=\"],\n    [\"61-35=\", \"61-4=\"],\n    [\"9+18=\", \"77-24=\"],\n    [\"10+30=\", \"94-73=\"],\n    [\"76+21=\", \"91-24=\"],\n    [\"55-28=\", \"71-31=\"],\n    [\"74+5=\", \"46-4=\"],\n    [\"34-16=\", \"14+1=\"],\n    [\"65-64=\", \"30+25=\"],\n    [\"75-8=\", \"34+38=\"],\n    [\"54+19=\", \"33+11=\"],\n    [\"88-82=\", \"89-83=\"],\n    [\"68+28=\", \"62-54=\"],\n    [\"19+35=\", \"5+28=\"],\n    [\"50-0=\", \"36+17=\"],\n    [\"45-42=\", \"82-26=\"],\n    [\"57-3=\", \"40-38=\"],\n    [\"81-53=\", \"99-65=\"],\n    [\"5+4=\", \"2+26=\"],\n    [\"57+17=\", \"5+62=\"],\n    [\"26-2=\", \"53-33=\"],\n    [\"11+87=\", \"50+34=\"],\n    [\"47-7=\", \"68+24=\"],\n    [\"60-31=\", \"3+39=\"],\n    [\"31+18=\", \"10+10=\"],\n    [\"21+27=\", \"41-8=\"],\n    [\"40+17=\", \"12+74=\"],\n    [\"48+44=\", \"11+67=\"],\n    [\"89-81=\", \"19+69=\"],\n    [\"99-58=\", \"19-8=\"],\n    [\"66-12=\", \"97-76=\"],\n    [\"73-23=\", \"8-4=\"],\n    [\"90-89=\", \"72-10=\"],\n    [\"18-10=\", \"59-35=\"],\n    [\"87-2=\", \"84+12=\"],\n    [\"55-38=\", \"56+6=\"],\n    [\"41+31=\", \"62+35=\"],\n    [\"59-51=\", \"63+3=\"],\n    [\"58-0=\", \"38+49=\"],\n    [\"65+24=\", \"24+42=\"],\n    [\"41+0=\", \"29-12=\"],\n    [\"17+18=\", \"62-24=\"],\n    [\"47+12=\", \"77-8=\"],\n    [\"31+11=\", \"6+67=\"],\n    [\"91-40=\", \"60+6=\"],\n    [\"64-15=\", \"27+18=\"],\n    [\"7+0=\", \"88-41=\"],\n    [\"66-6=\", \"31+4=\"],\n    [\"21+37=\", \"86-25=\"],\n    [\"79-15=\", \"37+62=\"],\n    [\"97-31=\", \"34-7=\"],\n    [\"45-19=\", \"98-9=\"],\n    [\"53-22=\", \"96-80=\"],\n    [\"83-47=\", \"47+39=\"],\n    [\"30+35=\", \"31-30=\"],\n    [\"66+32=\", \"43+18=\"],\n    [\"73-26=\", \"65+25=\"],\n    [\"62-45=\", \"58-22=\"],\n    [\"23+63=\", \"6+45=\"],\n    [\"85-26=\", \"75-23=\"],\n    [\"43+43=\", \"47-12=\"],\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load('items');\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error('Expected a table in the document, found none.');\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load('items');\nawait context.sync();\n\n// Load all cells' body paragraphs/text up front.\nconst cellParagraphsByRow = [];\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load('items');\n  cellParagraphsByRow.push(cells);\n}\nawait context.sync();\n\n// Flatten cells in row-major (document) order and load their text ranges.\nconst allCells = [];\nfor (const cells of cellParagraphsByRow) {\n  for (const cell of cells.items) {\n    allCells.push(cell);\n  }\n}\n\nconst ranges = allCells.map((cell) => {\n  const r = cell.body.getRange();\n  r.load('text');\n  return r;\n});\nawait context.sync();\n\nif (allCells.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} cells, found ${allCells.length}.`\n  );\n}\n\nfor (let i = 0; i < ranges.length; i++) {\n  const range = ranges[i];\n  const [oldText, newText] = replacements[i];\n  const actual = range.text.trim();\n  if (actual !== oldText) {\n    throw new Error(\n      `Cell ${i}: expected text \"${oldText}\" but found \"${actual}\".`\n    );\n  }\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replacement pairs extracted from the diff, in document order (old text -> new text).\n# All old texts are unique within the document, so this is unambiguous.\n$replacements = @(\n    @{ Old = '96-9='; New = '3+61=' },\n    @{ Old = '55-41='; New = '75-30=' },\n    @{ Old = '26+40='; New = '33+40=' },\n    @{ Old = '88-2='; New = '63-12=' },\n    @{ Old = '87-46='; New = '19+20=' },\n    @{ Old = '54-3='; New = '51+44=' },\n    @{ Old = '41+26='; New = '81-20=' },\n    @{ Old = '82-3='; New = '46-29=' },\n    @{ Old = '32+42='; New = '82-19=' },\n    @{ Old = '84+1='; New = '0+52=' },\n    @{ Old = '17+1='; New = '77-50=' },\n    @{ Old = '68-22='; New = '42+44=' },\n    @{ Old = '25+17='; New = '68-7=' },\n    @{ Old = '95-86='; New = '9+64=' },\n    @{ Old = '72-65='; New = '3+71=' },\n    @{ Old = '29+32='; New = '45-44=' },\n    @{ Old = '27+57='; New = '47+44=' },\n    @{ Old = '40+27='; New = '87-64=' },\n    @{ Old = '32+38='; New = '53+5=' },\n    @{ Old = '86-48='; New = '8+14=' },\n    @{ Old = '25+73='; New = '59+22=' },\n    @{ Old = '73-49='; New = '42-32=' },\n    @{ Old = '11+19='; New = '90-80=' },\n    @{ Old = '89-19='; New = '54+40=' },\n    @{ Old = '37-22='; New = '68-65=' },\n    @{ Old = '52-24='; New = '17+56=' },\n    @{ Old = '80-4='; New = '2+71=' },\n    @{ Old = '38+50='; New = '22-11=' },\n    @{ Old = '13+69='; New = '46+8=' },\n    @{ Old = '16+78='; New = '33+34=' },\n    @{ Old = '36+58='; New = '14+14=' },\n    @{ Old = '30+55='; New = '13+9=' },\n    @{ Old = '39-30='; New = '91-2=' },\n    @{ Old = '23-17='; New = '88-0=' },\n    @{ Old = '12+42='; New = '4+16=' },\n    @{ Old = '5+0='; New = '11+67=' },\n    @{ Old = '32+9='; New = '79-33=' },\n    @{ Old = '4+49='; New = '69+27=' },\n    @{ Old = '40-25='; New = '82-29=' },\n    @{ Old = '92-0='; New = '28+9=' },\n    @{ Old = '61-35='; New = '61-4=' },\n    @{ Old = '9+18='; New = '77-24=' },\n    @{ Old = '10+30='; New = '94-73=' },\n    @{ Old = '76+21='; New = '91-24=' },\n    @{ Old = '55-28='; New = '71-31=' },\n    @{ Old = '74+5='; New = '46-4=' },\n    @{ Old = '34-16='; New = '14+1=' },\n    @{ Old = '65-64='; New = '30+25=' },\n    @{ Old = '75-8='; New = '34+38=' },\n    @{ Old = '54+19='; New = '33+11=' },\n    @{ Old = '88-82='; New = '89-83=' },\n    @{ Old = '68+28='; New = '62-54=' },\n    @{ Old = '19+35='; New = '5+28=' },\n    @{ Old = '50-0='; New = '36+17=' },\n    @{ Old = '45-42='; New = '82-26=' },\n    @{ Old = '57-3='; New = '40-38=' },\n    @{ Old = '81-53='; New = '99-65=' },\n    @{ Old = '5+4='; New = '2+26=' },\n    @{ Old = '57+17='; New = '5+62=' },\n    @{ Old = '26-2='; New = '53-33=' },\n    @{ Old = '11+87='; New = '50+34=' },\n    @{ Old = '47-7='; New = '68+24=' },\n    @{ Old = '60-31='; New = '3+39=' },\n    @{ Old = '31+18='; New = '10+10=' },\n    @{ Old = '21+27='; New = '41-8=' },\n    @{ Old = '40+17='; New = '12+74=' },\n    @{ Old = '48+44='; New = '11+67=' },\n    @{ Old = '89-81='; New = '19+69=' },\n    @{ Old = '99-58='; New = '19-8=' },\n    @{ Old = '66-12='; New = '97-76=' },\n    @{ Old = '73-23='; New = '8-4=' },\n    @{ Old = '90-89='; New = '72-10=' },\n    @{ Old = '18-10='; New = '59-35=' },\n    @{ Old = '87-2='; New = '84+12=' },\n    @{ Old = '55-38='; New = '56+6=' },\n    @{ Old = '41+31='; New = '62+35=' },\n    @{ Old = '59-51='; New = '63+3=' },\n    @{ Old = '58-0='; New = '38+49=' },\n    @{ Old = '65+24='; New = '24+42=' },\n    @{ Old = '41+0='; New = '29-12=' },\n    @{ Old = '17+18='; New = '62-24=' },\n    @{ Old = '47+12='; New = '77-8=' },\n    @{ Old = '31+11='; New = '6+67=' },\n    @{ Old = '91-40='; New = '60+6=' },\n    @{ Old = '64-15='; New = '27+18=' },\n    @{ Old = '7+0='; New = '88-41=' },\n    @{ Old = '66-6='; New = '31+4=' },\n    @{ Old = '21+37='; New = '86-25=' },\n    @{ Old = '79-15='; New = '37+62=' },\n    @{ Old = '97-31='; New = '34-7=' },\n    @{ Old = '45-19='; New = '98-9=' },\n    @{ Old = '53-22='; New = '96-80=' },\n    @{ Old = '83-47='; New = '47+39=' },\n    @{ Old = '30+35='; New = '31-30=' },\n    @{ Old = '66+32='; New = '43+18=' },\n    @{ Old = '73-26='; New = '65+25=' },\n    @{ Old = '62-45='; New = '58-22=' },\n    @{ Old = '23+63='; New = '6+45=' },\n    @{ Old = '85-26='; New = '75-23=' },\n    @{ Old = '43+43='; New = '47-12=' }\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$expectedCount = $replacements.Count\n$actualCount = $rowCount * $colCount\nif ($actualCount -ne $expectedCount) {\n    throw \"Expected $expectedCount cells ($rowCount x $colCount = $actualCount), but counts differ.\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $pair = $replacements[$i]\n        $cell = $table.Cell($r, $c)\n        $rng = $cell.Range\n        # Range.Text for a table cell includes the trailing cell-mark characters (\\r\\a);\n        # trim them off before comparing against the expected value.\n        $current = $rng.Text.TrimEnd([char]13, [char]7)\n        if ($current -ne $pair.Old) {\n            throw \"Cell ($r,$c): expected '$($pair.Old)' but found '$current'.\"\n        }\n        $rng.Text = $pair.New\n        $i++\n    }\n}\n"}
